$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new rows of data to the worksheet
$ws.Range("A10").Value = "LookupValue"
$ws.Range("B10").Value = "Permissions"
$ws.Range("A11").Value = "LookupValue"
$ws.Range("B11").Value = "Copy Document"

# Resize the Excel Table (ListObject) to include the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B11"))

# Update the active selection to match the target state
$ws.Range("C16").Select()
